$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "initial_dipwell_measurements"
$ws.Range("B15").Value = "initial_condition/initial_day_dipwell_coords_and_measurements.csv"

$ws.Range("B15").Select()
